{"js": "// The document contains four transcription paragraphs whose \"<id>...</id>\"\n// marker was originally split across three runs:\n//   run1 (Courier New, color 7f6000): \"<id>\"\n//   run2 (default color):             \"<id_value>\"   (e.g. \"p164r_1\")\n//   run3 (Courier New, color 7f6000): \"</id>\"\n// The edit merges these three runs into a single run (keeping the Courier\n// New styling of the surrounding tag runs) whose text is the full\n// \"<id>...</id>\" string. One of the four ids (\"164r_4\") is also corrected\n// to \"p164r_4\" as part of the merge.\n\nconst idReplacements = [\n  { search: \"<id>p164r_1</id>\", replacement: \"<id>p164r_1</id>\" },\n  { search: \"<id>p164r_2</id>\", replacement: \"<id>p164r_2</id>\" },\n  { search: \"<id>p164r_3</id>\", replacement: \"<id>p164r_3</id>\" },\n  { search: \"<id>164r_4</id>\", replacement: \"<id>p164r_4</id>\" },\n];\n\nfor (const { search, replacement } of idReplacements) {\n  const results = context.document.body.search(search, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains four transcription paragraphs whose \"<id>...</id>\"\n# marker was originally split across three runs:\n#   run1 (Courier New, color 7f6000): \"<id>\"\n#   run2 (default color):             \"<id_value>\"   (e.g. \"p164r_1\")\n#   run3 (Courier New, color 7f6000): \"</id>\"\n# The edit merges these three runs into a single run (keeping the Courier\n# New styling of the surrounding tag runs) whose text is the full\n# \"<id>...</id>\" string. One of the four ids (\"164r_4\") is also corrected\n# to \"p164r_4\" as part of the merge.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Merge-IdTag($SearchText, $ReplaceText) {\n    $find = $word.ActiveDocument.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $SearchText\n    $find.Replacement.Text = $ReplaceText\n    $find.Execute($SearchText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $ReplaceText, $wdReplaceAll)\n}\n\nMerge-IdTag \"<id>p164r_1</id>\" \"<id>p164r_1</id>\"\nMerge-IdTag \"<id>p164r_2</id>\" \"<id>p164r_2</id>\"\nMerge-IdTag \"<id>p164r_3</id>\" \"<id>p164r_3</id>\"\nMerge-IdTag \"<id>164r_4</id>\" \"<id>p164r_4</id>\"\n"}
